$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-03-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-29 Saturday", 2) | Out-Null

# New answer values for the 20x5 answer table, in row-major (reading) order
$newValues = @(
    "45-8=37",
    "16+8=24",
    "69+19=88",
    "29+39=68",
    "19+28=47",
    "6+65=71",
    "36+39=75",
    "39+43=82",
    "36+55=91",
    "90-61=29",
    "87-58=29",
    "16+6=22",
    "84+9=93",
    "77+7=84",
    "20-16=4",
    "63-55=8",
    "13+78=91",
    "59+6=65",
    "8+74=82",
    "77+5=82",
    "37+59=96",
    "49+46=95",
    "28+43=71",
    "27+7=34",
    "19+19=38",
    "90-51=39",
    "60-46=14",
    "48-29=19",
    "89+9=98",
    "67-48=19",
    "30-19=11",
    "91-57=34",
    "46+15=61",
    "81-47=34",
    "52-47=5",
    "74-9=65",
    "17+47=64",
    "56-8=48",
    "6+58=64",
    "62+19=81",
    "64-8=56",
    "15+56=71",
    "93-37=56",
    "60-2=58",
    "83-78=5",
    "9+52=61",
    "47+17=64",
    "31-5=26",
    "65-56=9",
    "26+49=75",
    "18+77=95",
    "43+28=71",
    "51-8=43",
    "83-18=65",
    "62-48=14",
    "53-39=14",
    "58-29=29",
    "9+56=65",
    "75-69=6",
    "62-27=35",
    "85-6=79",
    "36-28=8",
    "91-58=33",
    "91-7=84",
    "19+37=56",
    "18+69=87",
    "39+25=64",
    "70-62=8",
    "26+25=51",
    "8+37=45",
    "97-38=59",
    "53-48=5",
    "74-47=27",
    "27+64=91",
    "9+27=36",
    "14+8=22",
    "53-24=29",
    "48+25=73",
    "38+47=85",
    "94-55=39",
    "69+4=73",
    "83-65=18",
    "83-58=25",
    "26+57=83",
    "63-15=48",
    "49+34=83",
    "39+47=86",
    "94-25=69",
    "26+69=95",
    "50-24=26",
    "18+17=35",
    "19+26=45",
    "83-79=4",
    "16+19=35",
    "92-36=56",
    "42-38=4",
    "40-24=16",
    "64-36=28",
    "95-8=87",
    "51-48=3"
)

$t = $d.Tables.Item(1)
$numRows = $t.Rows.Count
$numCols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "cells"